$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 714
    $ws.Range("F7").Value = 24
    $ws.Range("F12").Value = 4416
    $ws.Range("F15").Value = 152
}
